# Add a "2022-Q1" sheet (holdings detail) positioned between "2021-Q4" and
# "总计", and update the "总计" (totals) sheet with a new row summarizing
# 2022-Q1 while keeping the existing 2021-Q4 summary row.

$wb = $excel.ActiveWorkbook

$fmtFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

function Set-TextCell($cell, $val) {
    # Force the cell to stay text (Excel COM would otherwise infer a
    # number from digit-only / numeric-looking strings).
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet by duplicating "2021-Q4" (this keeps
#    sheetPr/pageMargins/column formatting identical) right after it,
#    then trim it down to the header + 2 data rows and overwrite values.
# ---------------------------------------------------------------------
$sheetQ4 = $wb.Worksheets.Item("2021-Q4")
$sheetQ4.Copy($null, $sheetQ4)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# The template sheet has 4 data rows (rows 2-5); 2022-Q1 only has 2, so
# drop the extra rows 4 and 5 (deleting row 4 three times removes both
# of the trailing rows after the first two).
$newSheet.Rows.Item(4).Delete()
$newSheet.Rows.Item(4).Delete()
$newSheet.Rows.Item(4).Delete()

# Row 2: 九泰久稳灵活配置混合A
Set-TextCell $newSheet.Cells.Item(2, 2) "002453"
Set-TextCell $newSheet.Cells.Item(2, 3) "九泰久稳灵活配置混合A"
Set-TextCell $newSheet.Cells.Item(2, 4) "0.09"
Set-TextCell $newSheet.Cells.Item(2, 5) "94.85"
Set-TextCell $newSheet.Cells.Item(2, 6) "1.92"
Set-TextCell $newSheet.Cells.Item(2, 7) "0.0017"
$newSheet.Cells.Item(2, 8).Value = 7

# Row 3: 九泰久稳灵活配置混合C
Set-TextCell $newSheet.Cells.Item(3, 2) "002454"
Set-TextCell $newSheet.Cells.Item(3, 3) "九泰久稳灵活配置混合C"
Set-TextCell $newSheet.Cells.Item(3, 4) "0.04"
Set-TextCell $newSheet.Cells.Item(3, 5) "94.85"
Set-TextCell $newSheet.Cells.Item(3, 6) "1.92"
Set-TextCell $newSheet.Cells.Item(3, 7) "0.0008"
$newSheet.Cells.Item(3, 8).Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: move the existing 2021-Q4 summary row from
#    row 2 down to row 3, then write the new 2022-Q1 summary into row 2.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$oldDate = $totalSheet.Cells.Item(2, 2).Value2
$oldCount = $totalSheet.Cells.Item(2, 3).Value2
$oldValue = $totalSheet.Cells.Item(2, 4).Value2

# Row 3 <- former row 2 (2021-Q4), including column A's style (s=2).
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial($fmtFormats)
$totalSheet.Cells.Item(3, 2).Value = $oldDate
$totalSheet.Cells.Item(3, 3).Value = $oldCount
$totalSheet.Cells.Item(3, 4).Value = $oldValue

# Row 2 <- new 2022-Q1 summary.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0

# Restore the original active sheet ("2021-Q4") so the workbook-level
# selection state isn't perturbed by all the sheet copying above.
$sheetQ4.Activate()

Write-Output "applied"
